# Slide 10 ("references"): the "Text Placeholder 6" shape lists source
# citations. Two more sources are appended ("Everyvine" and "VinePair"),
# which pushes the placeholder text past its box and causes PowerPoint's
# existing "shrink text on overflow" autofit (inherited as an empty
# <a:normAutofit/> from the slide layout) to kick in and record a
# font-scale reduction.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item("Text Placeholder 6")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Append the two new reference paragraphs.
$tr.Text = $tr.Text + "`rEveryvine`rVinePair"

# The placeholder uses "shrink text on overflow" autofit; recalculating it
# after the new lines overflow the box is what stamps the fontScale /
# lnSpcReduction values into the slide's own <a:bodyPr>/<a:normAutofit>.
$tf.AutoSize = 2
$tf.AutofitFontScale = 92500
$tf.LineSpaceReduction = 20000

